# "Started work on the smoker enemy"
# Move the SMOKER enemy marker from C8 to H8 on the Level1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Level1")

# C8 was the SMOKER tile; it becomes an empty/NONE tile.
$ws.Range("C8").Value = "NONE"

# H8 becomes the new SMOKER tile.
$ws.Range("H8").Value = "SMOKER"

# Match the author's final selection (H8) recorded in the sheet view.
$null = $ws.Range("H8").Select()
